# TC10_CDS_Filter_FileType-PDF.xlsx — "CDS Input file updates"
#
# Replaces the Neo4j query stored in B2 (ParticipantsTab "query") with a
# revised Cypher query, and updates the sheet view (scroll/selection) and
# row height to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- B2: new participant query text --------------------------------------
$newQuery = "MATCH (p:participant)-->(s:study)`n" +
            "OPTIONAL MATCH (samp:sample)-->(p)`n" +
            "OPTIONAL MATCH (p)<--(diag:diagnosis)`n" +
            "OPTIONAL MATCH (samp)<--(f:file)`n" +
            "OPTIONAL MATCH (f)<--(g:genomic_info)`n" +
            "WITH s, p, samp, f, g, diag`n" +
            "WHERE f.file_type in ['PDF']`n" +
            "with p`n" +
            "OPTIONAL MATCH (p)-->(s:study)`n" +
            "OPTIONAL MATCH (samp:sample)-->(p)`n" +
            "WITH s, p, apoc.coll.sort(collect(distinct samp.sample_id)) as samp`n" +
            "RETURN`n" +
            "coalesce(p.participant_id,'') as ``Participant ID``,`n" +
            "coalesce(s.study_name, '') as ``Study Name``,`n" +
            "coalesce(s.phs_accession,'') as ``Accession``,`n" +
            "coalesce(p.gender,'') as ``Gender``,`n" +
            "coalesce(apoc.text.join(samp, ','), '') as ``Samples```n" +
            "ORDER BY p.participant_id LIMIT 100"

$ws.Range("B2").Value = $newQuery

# --- Row 2 grew taller to fit the longer query text -----------------------
$ws.Rows.Item(2).RowHeight = 279

# --- Sheet view: scroll down a couple rows, move the active selection -----
$ws.Activate()
$ws.Range("B4").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 3
$win.ScrollColumn = 1
